$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new results row (row 3) to the quiz results sheet,
# mirroring the existing data rows (plain values, no special styling).
$ws.Cells.Item(3, 1).Value = 24152
$ws.Cells.Item(3, 2).Value = "kk"
$ws.Cells.Item(3, 3).Value = 9
$ws.Cells.Item(3, 4).Value = "2025-08-30 11:10:48"
